$d = $word.ActiveDocument

$d.Content.Find.Execute("808÷9=89, 7", $true, $false, $false, $false, $false, $true, 1, $false, "134÷8=16, 6", 2) | Out-Null
$d.Content.Find.Execute("777÷4=194, 1", $true, $false, $false, $false, $false, $true, 1, $false, "738÷7=105, 3", 2) | Out-Null
$d.Content.Find.Execute("948÷4=237, 0", $true, $false, $false, $false, $false, $true, 1, $false, "143÷7=20, 3", 2) | Out-Null
$d.Content.Find.Execute("110÷8=13, 6", $true, $false, $false, $false, $false, $true, 1, $false, "870÷2=435, 0", 2) | Out-Null
$d.Content.Find.Execute("108÷9=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "173÷9=19, 2", 2) | Out-Null
$d.Content.Find.Execute("876÷5=175, 1", $true, $false, $false, $false, $false, $true, 1, $false, "734÷3=244, 2", 2) | Out-Null
$d.Content.Find.Execute("113÷7=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "117÷4=29, 1", 2) | Out-Null
$d.Content.Find.Execute("314÷9=34, 8", $true, $false, $false, $false, $false, $true, 1, $false, "692÷4=173, 0", 2) | Out-Null
$d.Content.Find.Execute("981÷9=109, 0", $true, $false, $false, $false, $false, $true, 1, $false, "377÷3=125, 2", 2) | Out-Null
$d.Content.Find.Execute("551÷6=91, 5", $true, $false, $false, $false, $false, $true, 1, $false, "457÷9=50, 7", 2) | Out-Null
$d.Content.Find.Execute("921÷9=102, 3", $true, $false, $false, $false, $false, $true, 1, $false, "146÷7=20, 6", 2) | Out-Null
$d.Content.Find.Execute("273÷7=39, 0", $true, $false, $false, $false, $false, $true, 1, $false, "578÷2=289, 0", 2) | Out-Null
$d.Content.Find.Execute("218÷4=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "845÷7=120, 5", 2) | Out-Null
$d.Content.Find.Execute("552÷7=78, 6", $true, $false, $false, $false, $false, $true, 1, $false, "575÷8=71, 7", 2) | Out-Null
$d.Content.Find.Execute("307÷4=76, 3", $true, $false, $false, $false, $false, $true, 1, $false, "230÷2=115, 0", 2) | Out-Null
$d.Content.Find.Execute("460÷8=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "430÷3=143, 1", 2) | Out-Null
$d.Content.Find.Execute("615÷2=307, 1", $true, $false, $false, $false, $false, $true, 1, $false, "657÷3=219, 0", 2) | Out-Null
$d.Content.Find.Execute("931÷7=133, 0", $true, $false, $false, $false, $false, $true, 1, $false, "961÷9=106, 7", 2) | Out-Null
$d.Content.Find.Execute("319÷6=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "490÷2=245, 0", 2) | Out-Null
$d.Content.Find.Execute("191÷5=38, 1", $true, $false, $false, $false, $false, $true, 1, $false, "631÷6=105, 1", 2) | Out-Null
$d.Content.Find.Execute("378÷8=47, 2", $true, $false, $false, $false, $false, $true, 1, $false, "659÷7=94, 1", 2) | Out-Null
$d.Content.Find.Execute("124÷6=20, 4", $true, $false, $false, $false, $false, $true, 1, $false, "130÷2=65, 0", 2) | Out-Null
$d.Content.Find.Execute("816÷5=163, 1", $true, $false, $false, $false, $false, $true, 1, $false, "742÷2=371, 0", 2) | Out-Null
$d.Content.Find.Execute("426÷7=60, 6", $true, $false, $false, $false, $false, $true, 1, $false, "920÷4=230, 0", 2) | Out-Null
$d.Content.Find.Execute("622÷8=77, 6", $true, $false, $false, $false, $false, $true, 1, $false, "520÷7=74, 2", 2) | Out-Null
